$d = $word.ActiveDocument
$r = $d.Range($d.Content.End, $d.Content.End)
$xml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
    <w:p/>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">Upon further research, it was confirmed that there are </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">four different personality </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>types</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> and each type can be grouped into either passive or aggressive (</w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Barteau</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">). </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Barteau</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">wrote an article published in </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:i/>
          <w:iCs/>
        </w:rPr>
        <w:t xml:space="preserve">Dressage Today </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">in 2007. She is a U.S. national champion dressage rider, and discusses the four types of horse personalities that are seen in domestic horses. The types </w:t>
      </w:r>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>are;</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
      <w:r>
        <w:t xml:space="preserve"> social, fearful, challenging, and aloof. She goes into detail with the characteristics of each type and the “1-10” scale of those personalities. She also states that there is a passive to aggressive scale that applies to each type, with examples of a passive and aggressive version of each personality. Finally, the article goes into how to determine </w:t>
      </w:r>
      <w:r>
        <w:t>a specific</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> horse</w:t>
      </w:r>
      <w:r>
        <w:t>’</w:t>
      </w:r>
      <w:r>
        <w:t>s type and which behaviors and reactions can help you identify the type of personality you are looking at.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>B</w:t>
      </w:r>
      <w:r>
        <w:t>arteau</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>goes on in the article to give examples of how to determine which personality type specific horses are, which was utilized in this study.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t xml:space="preserve">An article by Foster also spoke to horse personalities, identifying the same four types as </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:t>Barteau</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">The main focus of Foster’s </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">article is </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">the way horses express discomfort with minimal movement. She identifies that </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">the eyes and other facial indicators are the most informative </w:t>
      </w:r>
      <w:r>
        <w:lastRenderedPageBreak/>
        <w:t>signals and that changes in body posture and natural movement are other signals to how a horse feels in turnout</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (2019).</w:t>
      </w:r>
    </w:p>

'@
$r.InsertXML($xml) | Out-Null
